$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Insert a brand-new paragraph right before "For each year, calculate..."
#    with a 0.5" (720 twips / 36 pt) left indent and the new commentary text.
# ---------------------------------------------------------------------------
$anchorPara = $d.Paragraphs.Item(2)
$anchorPara.Range.InsertParagraphBefore() | Out-Null

$newPara = $d.Paragraphs.Item(2)
$newPara.Range.Text = "For what the question asked, I didn" + [char]0x2019 + "t find any cases of the total volume crossing 1.5 billion for a given month of a year for any region. This was crosschecked in excel. incase this was meant to mean millions I have included the file " + [char]0x2018 + "millions.csv" + [char]0x2019 + " that contains the information in regards to 1.5 million total volume instead. "
$newPara.LeftIndent = 36

# ---------------------------------------------------------------------------
# 2) "The years most sold SKU by volume is as follows:" paragraph:
#    - "The" becomes "Each"
#    - a (now-empty) "_GoBack" bookmark is re-inserted between "vol" and "ume"
# ---------------------------------------------------------------------------
$skuPara = $d.Paragraphs.Item(4)
$skuRange = $skuPara.Range
$skuText = $skuRange.Text
$skuStart = $skuRange.Start

$theIdx = $skuText.IndexOf("The years")
$theStart = $skuStart + $theIdx
$theEnd = $theStart + 3
$d.Range($theStart, $theEnd).Text = "Each"

$skuPara = $d.Paragraphs.Item(4)
$skuRange = $skuPara.Range
$skuText = $skuRange.Text
$skuStart = $skuRange.Start
$volIdx = $skuText.IndexOf("volume")
$splitPos = $skuStart + $volIdx + 3
$splitRange = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("_GoBack", $splitRange) | Out-Null

# ---------------------------------------------------------------------------
# 3) Wrap "TotalUS" in spell-check proof-error markers inside the
#    "For each category..." paragraph. We rebuild the whole paragraph via
#    InsertXML so the markers land inline (not appended out of place).
# ---------------------------------------------------------------------------
$totalUsPara = $d.Paragraphs.Item(9)
$paraXmlNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$innerPara = '<w:p ' + $paraXmlNs + '><w:r><w:t xml:space="preserve">For each category and each month, calculate the ratio of large bags sold in the &#8220;</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>TotalUS</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>&#8221; region to large bags sold in the region &#8220;Charlotte.&#8221;</w:t></w:r></w:p>'
$packageXml = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $innerPara + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$totalUsPara.Range.InsertXML($packageXml)

# ---------------------------------------------------------------------------
# 4) Remove the old "_GoBack" bookmark from the final paragraph (it moved to
#    the SKU paragraph above).
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $lastParaText = $d.Paragraphs.Item($d.Paragraphs.Count).Range.Text
    if ($lastParaText -like "*data is arranged*") {
        $bm = $d.Bookmarks.Item("_GoBack")
        $bmRange = $bm.Range
        if ($bmRange.Start -gt $skuRange.End) {
            $bm.Delete()
        }
    }
}

Write-Output "edit complete"
